# Automatische test-sync: 2025-08-13 20:58:50
# Append a new log entry row to the "Logs" sheet, extend the conditional
# formatting ranges to cover it, and refresh the "Dashboard" sheet's
# summary count.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 9

$logs.Cells.Item($newRow, 1).Value = "Demo inplannen"
$logs.Cells.Item($newRow, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item($newRow, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-13 20:58:14"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the existing conditional-formatting rules (columns D, G, H, I, J)
# so they cover the newly added row, keeping the same rules/dxf mapping.
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "8")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "9")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 8
